$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-9: set Absent (H) to 1
foreach ($r in 3..9) {
    $ws.Range("H$r").Value = 1
}

# Rows 10-15: set Total Attendance Count (D) and Real (E) to 1
foreach ($r in 10..15) {
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = 1
}

# Rows 16-18: set Absent (H) to 1
foreach ($r in 16..18) {
    $ws.Range("H$r").Value = 1
}
